$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.420.01"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "'2.066.66"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'234.81"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'56.98"
$ws.Range("E8").Value = "  -2.55%  "

$ws.Range("E9").Value = "  +2.21%  "

$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").Value = "'2.372.48"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "'14.30"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("D14").Value = "'20.59"
$ws.Range("E14").Value = "  -3.80%  "

$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "'2.067.95"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "'37.314.57"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  +1.02%  "

$ws.Range("D20").Value = "'69.50"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("D21").Value = "'0.0₃0819"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").Value = "'226.18"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("E25").Value = "  -2.51%  "

$ws.Range("D26").Value = "'167.92"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("D27").Value = "'8.84"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").Value = "'0.133"
$ws.Range("E28").Value = "  +4.69%  "

$ws.Range("E29").Value = "  -6.04%  "

$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("E33").Value = "  -1.12%  "

$ws.Range("D34").Value = "'4.53"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("D36").Value = "'3.36"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").Value = "'5.62"
$ws.Range("E39").Value = "  -4.32%  "

$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("D41").Value = "'1.492.16"
$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("D42").Value = "'0.0954"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("D43").Value = "'96.81"
$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").Value = "'4.20"
$ws.Range("E46").Value = "  -5.89%  "

$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("E48").Value = "  -4.16%  "

$ws.Range("D49").Value = "'7.21"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  +0.82%  "

$ws.Range("D51").Value = "'2.258.44"
$ws.Range("E51").Value = "  +0.00%  "

